# cadastro_monitores.xlsx - update monitor names and drop the extra rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old monitor names with the new list
$ws.Range("A2").Value = "julia"
$ws.Range("A3").Value = "maria"
$ws.Range("A4").Value = "ana"

# Remove the now-unused trailing rows (previously Louis, Zayn, Julia K)
$ws.Range("A5:A7").EntireRow.Delete()

# Reset selection back to the top of the sheet
$ws.Range("A1").Select()
